$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '42.840.56'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.299.48'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -1.49%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.119'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '2.659.89'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '2.297.16'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.782'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '42.781.19'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '0.0₃0900'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("E24").Value = '  -2.93%  '
$ws.Range("E25").Value = '  +1.83%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.12%  '
$ws.Range("E30").Value = '  +1.52%  '
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.87%  '
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").Value = '2.003.55'
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("E44").Value = '  -2.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.84%  '
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.13%  '
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = '2.528.37'
$ws.Range("E51").Value = '  -0.33%  '
